# Updates outdated/specific technology names in the "94Cram" marketing
# deck with more generic/"impressive" marketing language.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 3 - "AI 驅動核心" banner (single run)
# ---------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(30).TextFrame.TextRange.Text = `
  "🤖  AI 驅動核心：頂尖大型語言模型 + 智慧知識引擎 + 自然語言操作 + 智慧預測"

# ---------------------------------------------------------------
# Slide 10 - architecture diagram labels
# ---------------------------------------------------------------
$s10 = $p.Slides.Item(10)

# 前端層 (single run)
$s10.Shapes.Item(12).TextFrame.TextRange.Text = `
  "新一代響應式框架 · 強型別全棧開發 · 伺服器端渲染加速"

# API 層 (single run)
$s10.Shapes.Item(16).TextFrame.TextRange.Text = `
  "軍規級身份認證 · 角色權限管控 · 資料驗證防護 · 標準化介面"

# AI 層 (single run)
$s10.Shapes.Item(20).TextFrame.TextRange.Text = `
  "最新一代大型語言模型 · 智慧知識引擎 · 語意向量檢索 · 意圖理解"

# 資料層 (single run)
$s10.Shapes.Item(24).TextFrame.TextRange.Text = `
  "企業級關聯式資料庫 · 型別安全 ORM · 多租戶隔離 · 審計日誌"

# 基礎設施 (single run)
$s10.Shapes.Item(28).TextFrame.TextRange.Text = `
  "頂級雲端無伺服器架構 · 託管式資料庫 · 自動擴縮 · 零停機部署"

# "Blue-Green 部署" card — only the first run changes, the line
# break + second run ("更新系統完全不影響使用") stay untouched.
$tr32 = $s10.Shapes.Item(32).TextFrame.TextRange
$tr32.Characters(1, 13).Text = "雙軌熱切換部署"

# "JWT + RBAC + 審計日誌" card — only the first run changes.
$tr36 = $s10.Shapes.Item(36).TextFrame.TextRange
$tr36.Characters(1, 17).Text = "軍規級認證 + 角色權限"

# "Google Cloud 台灣區" card — only the first run changes.
$tr44 = $s10.Shapes.Item(44).TextFrame.TextRange
$tr44.Characters(1, 16).Text = "國際頂級雲端台灣區"

# ---------------------------------------------------------------
# Slide 11 - security & compliance grid
# ---------------------------------------------------------------
$s11 = $p.Slides.Item(11)

$s11.Shapes.Item(7).TextFrame.TextRange.Text = "軍規級認證 + SSO"
$s11.Shapes.Item(23).TextFrame.TextRange.Text = "頂級雲端防護"
